# Update the cryptocurrency price/volume table with the latest scrape.
# Column D ("Price") values are free-form numeric-looking strings (e.g.
# thousands-dot-separated "58.938.67", or tiny decimals) that must stay
# TEXT, so NumberFormat is forced to "@" (Text) before each write to stop
# Excel from auto-coercing them into floating point numbers and mangling
# the formatting (trailing zeros, multi-dot groupings, etc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.938.67'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.606.16'
$ws.Range('E3').Value = '  -1.24%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.36'
$ws.Range('E5').Value = '  +3.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.95'
$ws.Range('E6').Value = '  -1.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +4.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.83'
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.142'
$ws.Range('E11').Value = '  +5.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.336'
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.064.80'
$ws.Range('E13').Value = '  -1.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '58.902.98'
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.90'
$ws.Range('E15').Value = '  -1.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.614.08'
$ws.Range('E16').Value = '  -0.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('E17').Value = '  -1.70%  '
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '338.24'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.10'
$ws.Range('E20').Value = '  -1.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.67'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.428'
$ws.Range('E24').Value = '  +2.81%  '
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('E27').Value = '  -1.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0761'
$ws.Range('E28').Value = '  +1.86%  '
$ws.Range('E30').Value = '  +1.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.02'
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '154.26'
$ws.Range('E32').Value = '  +2.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.02'
$ws.Range('E33').Value = '  +1.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.97'
$ws.Range('E34').Value = '  -0.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.906'
$ws.Range('E35').Value = '  +8.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.895'
$ws.Range('E36').Value = '  +7.30%  '
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.09'
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '283.40'
$ws.Range('E41').Value = '  -0.53%  '
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0954'
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0537'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.62'
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.66'
$ws.Range('E48').Value = '  +2.41%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.955.19'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '119.03'
$ws.Range('E50').Value = '  +7.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.00'
$ws.Range('E51').Value = '  -2.03%  '
